# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# header style used by the other header cells (e.g. column F "Win").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 — same style as the other header cells in row 1
# (copy formatting from the neighboring "Win" header, then set the text).
$ws.Range("F1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for H2:H7.
$values = @(0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
